$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 34, shifting the existing rows 34-63 down to 35-64.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new weekly record, matching the
# structure/pattern of the surrounding rows (same market/category/etc.).
$ws.Range("A34").Value2 = 2
$ws.Range("B34").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C34").Value2 = "Coquimbo"
$ws.Range("D34").Value2 = 44763
$ws.Range("E34").Value2 = 4
$ws.Range("F34").Value2 = 100112022
$ws.Range("G34").Value2 = "Arveja Verde"
$ws.Range("H34").Value2 = "Perfection"
$ws.Range("I34").Value2 = "Primera"
$ws.Range("J34").Value2 = 400
$ws.Range("K34").Value2 = 29000
$ws.Range("L34").Value2 = 30000
$ws.Range("M34").Value2 = 29500
$ws.Range("N34").Value2 = "$/malla 25 kilos"
$ws.Range("O34").Value2 = "Provincia de Limarí"
$ws.Range("P34").Value2 = 1180
$ws.Range("Q34").Value2 = 25
$ws.Range("R34").Value2 = "Hortaliza"

# Make sure the new date cell keeps the same date number format as the rest
# of column D.
$ws.Range("D34").NumberFormat = $ws.Range("D35").NumberFormat
